# Auto-generated edit script: applies scheduled-runner market-data refresh
# to the Marilith_Profits workbook (updates computed price/profit columns H-N
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 66.14286
$ws.Range("J9").Value = 85
$ws.Range("L9").Value = 85
$ws.Range("N9").Value = -423
$ws.Range("H15").Value = 1848.12
$ws.Range("I15").Value = 1848.12
$ws.Range("K15").Value = 5544.36
$ws.Range("M15").Value = -5375.36
$ws.Range("H39").Value = 383.4
$ws.Range("J39").Value = 1000
$ws.Range("L39").Value = 3000
$ws.Range("N39").Value = -3592
$ws.Range("H42").Value = 193.2
$ws.Range("I42").Value = 216.5
$ws.Range("J42").Value = 100
$ws.Range("K42").Value = 649.5
$ws.Range("L42").Value = 300
$ws.Range("M42").Value = -419.5
$ws.Range("N42").Value = -760
$ws.Range("H82").Value = 1266
$ws.Range("I82").Value = 560.3333
$ws.Range("K82").Value = 1680.9999
$ws.Range("M82").Value = -1274.9999
$ws.Range("H85").Value = 1266
$ws.Range("I85").Value = 560.3333
$ws.Range("K85").Value = 1680.9999
$ws.Range("M85").Value = -276.9999
$ws.Range("H104").Value = 250
$ws.Range("I104").Value = 250
$ws.Range("K104").Value = 750
$ws.Range("M104").Value = 997
$ws.Range("H112").Value = 2191
$ws.Range("J112").Value = 2228.8462
$ws.Range("L112").Value = 6686.5386
$ws.Range("N112").Value = -8902.5386
$ws.Range("H113").Value = 4443.3
$ws.Range("I113").Value = 4179.8
$ws.Range("J113").Value = 4706.8
$ws.Range("K113").Value = 4179.8
$ws.Range("L113").Value = 4706.8
$ws.Range("M113").Value = -925.8000000000002
$ws.Range("N113").Value = -11214.8
$ws.Range("H132").Value = 7380.625
$ws.Range("I132").Value = 7579.2856
$ws.Range("K132").Value = 22737.8568
$ws.Range("M132").Value = -20207.8568
$ws.Range("H141").Value = 3353.4167
$ws.Range("I141").Value = 2749.4546
$ws.Range("K141").Value = 8248.363799999999
$ws.Range("M141").Value = -3068.363799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2083.238
$ws.Range("J88").Value = 2842.5715
$ws.Range("L88").Value = 2842.5715
$ws.Range("N88").Value = -3654.5715
$ws.Range("H91").Value = 2083.238
$ws.Range("J91").Value = 2842.5715
$ws.Range("L91").Value = 2842.5715
$ws.Range("N91").Value = -5650.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1883.5416
$ws.Range("I86").Value = 1182.7222
$ws.Range("J86").Value = 3986
$ws.Range("K86").Value = 1182.7222
$ws.Range("L86").Value = 3986
$ws.Range("M86").Value = -59.72219999999993
$ws.Range("N86").Value = -6232
$ws.Range("H89").Value = 1883.5416
$ws.Range("I89").Value = 1182.7222
$ws.Range("J89").Value = 3986
$ws.Range("K89").Value = 5913.611
$ws.Range("L89").Value = 19930
$ws.Range("M89").Value = -297.6109999999999
$ws.Range("N89").Value = -31162

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2485.2307
$ws.Range("I31").Value = 1409.75
$ws.Range("J31").Value = 4206
$ws.Range("K31").Value = 1409.75
$ws.Range("L31").Value = 4206
$ws.Range("M31").Value = -1114.75
$ws.Range("N31").Value = -4796
$ws.Range("H34").Value = 2485.2307
$ws.Range("I34").Value = 1409.75
$ws.Range("J34").Value = 4206
$ws.Range("K34").Value = 1409.75
$ws.Range("L34").Value = 4206
$ws.Range("M34").Value = -1207.75
$ws.Range("N34").Value = -4610
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970
$ws.Range("H134").Value = 7666.5557
$ws.Range("I134").Value = 6749.875
$ws.Range("K134").Value = 20249.625
$ws.Range("M134").Value = -17714.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 795.34283
$ws.Range("I132").Value = 775.4483
$ws.Range("J132").Value = 891.5
$ws.Range("K132").Value = 6979.0347
$ws.Range("L132").Value = 8023.5
$ws.Range("M132").Value = -4449.0347
$ws.Range("N132").Value = -13083.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2283.2
$ws.Range("I70").Value = 2253.5
$ws.Range("K70").Value = 2253.5
$ws.Range("M70").Value = -1983.5
$ws.Range("H73").Value = 2283.2
$ws.Range("I73").Value = 2253.5
$ws.Range("K73").Value = 2253.5
$ws.Range("M73").Value = -1317.5
$ws.Range("H113").Value = 1949.75
$ws.Range("J113").Value = 1999.6666
$ws.Range("L113").Value = 1999.6666
$ws.Range("N113").Value = -6339.6666
$ws.Range("H122").Value = 5210845
$ws.Range("I122").Value = 5684276.5
$ws.Range("K122").Value = 17052829.5
$ws.Range("M122").Value = -17050379.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4427.857
$ws.Range("I7").Value = 2997.5
$ws.Range("K7").Value = 2997.5
$ws.Range("M7").Value = -2885.5
$ws.Range("H22").Value = 1379.7
$ws.Range("I22").Value = 1537.5
$ws.Range("J22").Value = 748.5
$ws.Range("K22").Value = 1537.5
$ws.Range("L22").Value = 748.5
$ws.Range("M22").Value = -1242.5
$ws.Range("N22").Value = -1338.5
$ws.Range("H27").Value = 1379.7
$ws.Range("I27").Value = 1537.5
$ws.Range("J27").Value = 748.5
$ws.Range("K27").Value = 1537.5
$ws.Range("L27").Value = 748.5
$ws.Range("M27").Value = -1430.5
$ws.Range("N27").Value = -962.5
$ws.Range("H40").Value = 8330
$ws.Range("I40").Value = 7495
$ws.Range("K40").Value = 7495
$ws.Range("M40").Value = -7359
$ws.Range("H61").Value = 4806.6924
$ws.Range("I61").Value = 4790.1665
$ws.Range("J61").Value = 5005
$ws.Range("K61").Value = 4790.1665
$ws.Range("L61").Value = 5005
$ws.Range("M61").Value = -4588.1665
$ws.Range("N61").Value = -5409
$ws.Range("H113").Value = 4806.6924
$ws.Range("I113").Value = 4790.1665
$ws.Range("J113").Value = 5005
$ws.Range("K113").Value = 4790.1665
$ws.Range("L113").Value = 5005
$ws.Range("M113").Value = -2620.1665
$ws.Range("N113").Value = -9345
$ws.Range("H126").Value = 4427.857
$ws.Range("I126").Value = 2997.5
$ws.Range("K126").Value = 8992.5
$ws.Range("M126").Value = -6522.5
$ws.Range("H132").Value = 11530.6
$ws.Range("I132").Value = 11530.6
$ws.Range("K132").Value = 34591.8
$ws.Range("M132").Value = -32061.8
$ws.Range("H136").Value = 3503.8333
$ws.Range("I136").Value = 3503.8333
$ws.Range("K136").Value = 10511.4999
$ws.Range("M136").Value = -7961.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 12000
$ws.Range("I37").Value = 12000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 12000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -11797
$ws.Range("N37").ClearContents()
$ws.Range("H81").Value = 812.5
$ws.Range("J81").Value = 833.3333
$ws.Range("L81").Value = 1666.6666
$ws.Range("N81").Value = -3788.6666
$ws.Range("H84").Value = 812.5
$ws.Range("J84").Value = 833.3333
$ws.Range("L84").Value = 8333.333000000001
$ws.Range("N84").Value = -18941.333
$ws.Range("H122").Value = 2119.1
$ws.Range("I122").Value = 2119.1
$ws.Range("K122").Value = 6357.299999999999
$ws.Range("M122").Value = -3907.299999999999
